$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows 10-13 (data no longer present in the updated TPM output)
$ws.Rows("10:13").Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.427849
$ws.Range("H2").Value = 4.283547
$ws.Range("I2").Value = 0.1372193253303967
$ws.Range("J2").Value = 0.1372193253303967
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1797713333333334
$ws.Range("N2").Value = 0.5393140000000001
$ws.Range("O2").Value = 0.188800001120238
$ws.Range("P2").Value = 0.188800001120238
$ws.Range("Q2").Value = 0.2566863185286667
$ws.Range("R2").Value = 2.310176866758001
$ws.Range("S2").Value = 0.0259070087760972
$ws.Range("T2").Value = 0.0259070087760972

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.427849
$ws.Range("H3").Value = 4.283547
$ws.Range("I3").Value = 0.1372193253303967
$ws.Range("J3").Value = 0.1372193253303967
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7724073333333333
$ws.Range("N3").Value = 2.317222
$ws.Range("O3").Value = 0.811199998879762
$ws.Range("P3").Value = 0.811199998879762
$ws.Range("Q3").Value = 1.102881038492667
$ws.Range("R3").Value = 9.925929346434001
$ws.Range("S3").Value = 0.1113123165542995
$ws.Range("T3").Value = 0.1113123165542995

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.757808
$ws.Range("H4").Value = 11.273424
$ws.Range("I4").Value = 0.361133340066889
$ws.Range("J4").Value = 0.361133340066889
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1797713333333334
$ws.Range("N4").Value = 0.5393140000000001
$ws.Range("O4").Value = 0.188800001120238
$ws.Range("P4").Value = 0.188800001120238
$ws.Range("Q4").Value = 0.6755461545706668
$ws.Range("R4").Value = 6.079915391136001
$ws.Range("S4").Value = 0.06818197500918392
$ws.Range("T4").Value = 0.06818197500918392

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.757808
$ws.Range("H5").Value = 11.273424
$ws.Range("I5").Value = 0.361133340066889
$ws.Range("J5").Value = 0.361133340066889
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7724073333333333
$ws.Range("N5").Value = 2.317222
$ws.Range("O5").Value = 0.811199998879762
$ws.Range("P5").Value = 0.811199998879762
$ws.Range("Q5").Value = 2.902558456458667
$ws.Range("R5").Value = 26.123026108128
$ws.Range("S5").Value = 0.2929513650577051
$ws.Range("T5").Value = 0.2929513650577051

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.167149333333334
$ws.Range("H6").Value = 15.501448
$ws.Range("I6").Value = 0.4965740392726466
$ws.Range("J6").Value = 0.4965740392726465
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1797713333333334
$ws.Range("N6").Value = 0.5393140000000001
$ws.Range("O6").Value = 0.188800001120238
$ws.Range("P6").Value = 0.188800001120238
$ws.Range("Q6").Value = 0.928905325185778
$ws.Range("R6").Value = 8.360147926672001
$ws.Range("S6").Value = 0.09375317917095677
$ws.Range("T6").Value = 0.09375317917095675

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.167149333333334
$ws.Range("H7").Value = 15.501448
$ws.Range("I7").Value = 0.4965740392726466
$ws.Range("J7").Value = 0.4965740392726465
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7724073333333333
$ws.Range("N7").Value = 2.317222
$ws.Range("O7").Value = 0.811199998879762
$ws.Range("P7").Value = 0.811199998879762
$ws.Range("Q7").Value = 3.991144037495111
$ws.Range("R7").Value = 35.920296337456
$ws.Range("S7").Value = 0.4028208601016898
$ws.Range("T7").Value = 0.4028208601016898

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.05279066666666667
$ws.Range("H8").Value = 0.158372
$ws.Range("I8").Value = 0.00507329533006772
$ws.Range("J8").Value = 0.005073295330067719
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1797713333333334
$ws.Range("N8").Value = 0.5393140000000001
$ws.Range("O8").Value = 0.188800001120238
$ws.Range("P8").Value = 0.188800001120238
$ws.Range("Q8").Value = 0.009490248534222225
$ws.Range("R8").Value = 0.08541223680800002
$ws.Range("S8").Value = 0.0009578381640000835
$ws.Range("T8").Value = 0.0009578381640000833

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.05279066666666667
$ws.Range("H9").Value = 0.158372
$ws.Range("I9").Value = 0.00507329533006772
$ws.Range("J9").Value = 0.005073295330067719
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7724073333333333
$ws.Range("N9").Value = 2.317222
$ws.Range("O9").Value = 0.811199998879762
$ws.Range("P9").Value = 0.811199998879762
$ws.Range("Q9").Value = 0.04077589806488889
$ws.Range("R9").Value = 0.366983082584
$ws.Range("S9").Value = 0.004115457166067637
$ws.Range("T9").Value = 0.004115457166067636
